$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cells whose new value would be auto-parsed as a number by Excel (force text) ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D5").Value = "575.45"
$ws.Range("D6").Value = "167.31"
$ws.Range("D9").Value = "0.119"
$ws.Range("D11").Value = "0.392"
$ws.Range("D15").Value = "25.65"
$ws.Range("D18").Value = "412.75"
$ws.Range("D19").Value = "12.91"
$ws.Range("D23").Value = "69.75"
$ws.Range("D26").Value = "0.0000105"
$ws.Range("D27").Value = "8.92"
$ws.Range("D28").Value = "1.00"
$ws.Range("D29").Value = "1.85"
$ws.Range("D30").Value = "21.65"
$ws.Range("D31").Value = "4.99"
$ws.Range("D32").Value = "6.41"
$ws.Range("D34").Value = "157.34"
$ws.Range("D35").Value = "1.37"
$ws.Range("D38").Value = "24.31"
$ws.Range("D41").Value = "0.0634"
$ws.Range("D42").Value = "5.63"
$ws.Range("D44").Value = "296.90"
$ws.Range("D45").Value = "21.54"
$ws.Range("D47").Value = "0.0994"
$ws.Range("D48").Value = "1.97"
$ws.Range("D50").Value = "10.48"
$ws.Range("D51").Value = "0.911"

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

# --- Set cells whose new value already reads as text (has letters/percent/multi-dot/spaces) ---
$ws.Range("D2").Value = "65.116.41"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "3.212.44"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  -5.44%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "3.772.58"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "65.147.32"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "3.210.68"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("E26").Value = "  -5.60%  "
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "2.749.15"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").Value = "  -8.96%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("E51").Value = "  -2.05%  "
